# The source workbook tracks daily "Lechuga" (lettuce) price records for the
# "Feria Lagunitas de Puerto Montt" market. This edit inserts one new daily
# record as a new row 907 (pushing the existing row 907 and everything below
# it down by one row), matching the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 907; this shifts the old rows 907-946 down to 908-947.
$ws.Rows("907:907").Insert()

# Populate the newly inserted row 907 with the new record. Columns A, B, C,
# E, F, G, H, I, N, O, Q, R repeat the same values as the (now shifted) row
# below it; D, J, K, L, M, P carry the new record's own data.
$ws.Range("A907").Value = 4
$ws.Range("B907").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C907").Value = "Los Lagos"
$ws.Range("D907").Value = 45147
$ws.Range("D907").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E907").Value = 10
$ws.Range("F907").Value = 100112033
$ws.Range("G907").Value = "Lechuga"
$ws.Range("H907").Value = "Escarola"
$ws.Range("I907").Value = "Primera"
$ws.Range("J907").Value = 150
$ws.Range("K907").Value = 13000
$ws.Range("L907").Value = 13000
$ws.Range("M907").Value = 13000
$ws.Range("N907").Value = "$/caja 15 unidades"
$ws.Range("O907").Value = "Región de Coquimbo"
$ws.Range("P907").Value = 867
$ws.Range("Q907").Value = 15
$ws.Range("R907").Value = "Hortaliza"
